# feat: add 2022-Q1 data
#
# The previously-last sheet "总计" (totals) is renamed to "2022-Q1" and its
# contents are replaced with the per-fund holdings for that quarter (same
# layout as the other quarterly sheets, e.g. "2021-Q4"). A brand-new "总计"
# sheet is then appended at the end, containing the same totals table as
# before plus a new leading row for "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Locate the existing "总计" sheet and the "2021-Q4" sheet (used as a
#    formatting template, since it already has the fund-holdings layout).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 2. Rename the old "总计" sheet to "2022-Q1" and wipe its contents so we
#    can write the fresh fund-holdings table into it.
# ---------------------------------------------------------------------
$q1Sheet = $totalSheet
$q1Sheet.Name = "2022-Q1"
$q1Sheet.Cells.Clear()

# Copy header-row / index-column formatting from the template sheet so the
# new sheet matches the look of the other quarterly sheets.
$templateSheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$q1Sheet.Range("A2:A5").PasteSpecial(-4122)

$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Index column (A) and rank column (H) are genuine numbers; the rest of
# the row is kept as text (matches the source data, which keeps fund
# codes / percentages as plain strings e.g. to preserve leading zeros).
for ($i = 0; $i -le 3; $i++) {
    $q1Sheet.Cells.Item(2 + $i, 1).Value = $i
}

$textCols = $q1Sheet.Range("B2:G5")
$textCols.NumberFormat = "@"

$q1Sheet.Range("B2").Value = "010714"
$q1Sheet.Range("C2").Value = "东方红远见价值混合"
$q1Sheet.Range("D2").Value = "19.70"
$q1Sheet.Range("E2").Value = "86.34"
$q1Sheet.Range("F2").Value = "3.02"
$q1Sheet.Range("G2").Value = "0.5949"
$q1Sheet.Range("H2").Value = 8

$q1Sheet.Range("B3").Value = "011651"
$q1Sheet.Range("C3").Value = "招商港股通核心精选股票A"
$q1Sheet.Range("D3").Value = "2.81"
$q1Sheet.Range("E3").Value = "81.27"
$q1Sheet.Range("F3").Value = "2.69"
$q1Sheet.Range("G3").Value = "0.0756"
$q1Sheet.Range("H3").Value = 8

$q1Sheet.Range("B4").Value = "519139"
$q1Sheet.Range("C4").Value = "海富通沪港深灵活配置混合"
$q1Sheet.Range("D4").Value = "1.32"
$q1Sheet.Range("E4").Value = "94.37"
$q1Sheet.Range("F4").Value = "3.05"
$q1Sheet.Range("G4").Value = "0.0403"
$q1Sheet.Range("H4").Value = 10

$q1Sheet.Range("B5").Value = "011652"
$q1Sheet.Range("C5").Value = "招商港股通核心精选股票C"
$q1Sheet.Range("D5").Value = "0.94"
$q1Sheet.Range("E5").Value = "81.27"
$q1Sheet.Range("F5").Value = "2.69"
$q1Sheet.Range("G5").Value = "0.0253"
$q1Sheet.Range("H5").Value = 8

# Drop the explicit Text number format again so the cells end up with no
# special style applied, matching the rest of the workbook's data rows.
$textCols.ClearFormats()

# ---------------------------------------------------------------------
# 3. Add a brand-new "总计" sheet after "2022-Q1" containing the updated
#    totals table (previous rows plus a new "2022-Q1" row at the top).
# ---------------------------------------------------------------------
$newTotalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$newTotalSheet.Name = "总计"

$templateSheet.Range("B1:D1").Copy()
$newTotalSheet.Range("B1:D1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newTotalSheet.Range("A2:A7").PasteSpecial(-4122)

$newTotalSheet.Range("B1").Value = "日期"
$newTotalSheet.Range("C1").Value = "持有数量(只)"
$newTotalSheet.Range("D1").Value = "持有市值(亿元)"

for ($i = 0; $i -le 5; $i++) {
    $newTotalSheet.Cells.Item(2 + $i, 1).Value = $i
}

$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 4
$newTotalSheet.Range("D2").Value = 0.74

$newTotalSheet.Range("B3").Value = "2021-Q4"
$newTotalSheet.Range("C3").Value = 3
$newTotalSheet.Range("D3").Value = 1.15

$newTotalSheet.Range("B4").Value = "2021-Q3"
$newTotalSheet.Range("C4").Value = 1
$newTotalSheet.Range("D4").Value = 1

$newTotalSheet.Range("B5").Value = "2021-Q2"
$newTotalSheet.Range("C5").Value = 1
$newTotalSheet.Range("D5").Value = 1.14

$newTotalSheet.Range("B6").Value = "2021-Q1"
$newTotalSheet.Range("C6").Value = 2
$newTotalSheet.Range("D6").Value = 0.97

$newTotalSheet.Range("B7").Value = "2020-Q4"
$newTotalSheet.Range("C7").Value = 3
$newTotalSheet.Range("D7").Value = 1.7
